# Update cryptos list data (Price and Volume(1h) columns) per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.141.69"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.905.42"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "325.89"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "0.4610"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "0.3891"
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("D9").Value = "0.07879"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").Value = "0.9901"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").Value = "21.99"
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("D12").Value = "1.883.24"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "5.767"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "7.040"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "0.07036"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "88.12"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "29.180.31"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "5.324"
$ws.Range("D23").Value = "11.13"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "2.103"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").Value = "156.28"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").Value = "19.46"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("D28").Value = "118.85"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "1.879"
$ws.Range("E29").Value = "  -6.16%  "
$ws.Range("D30").Value = "0.09356"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").Value = "0.8959"
$ws.Range("E31").Value = "  -3.51%  "
$ws.Range("D32").Value = "5.229"
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("D33").Value = "1.322"
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("D34").Value = "3.157"
$ws.Range("E34").Value = "  -3.64%  "
$ws.Range("D35").Value = "0.05795"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").Value = "1.173"
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("D37").Value = "0.02088"
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("D38").Value = "1.001"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").Value = "7.681"
$ws.Range("E40").Value = "  -3.60%  "
$ws.Range("D41").Value = "0.1810"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").Value = "9.718"
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("D43").Value = "11.92"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").Value = "0.5359"
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").Value = "2.179"
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("D46").Value = "0.07010"
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").Value = "1.842"
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").Value = "113.15"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "0.2959"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").Value = "71.27"
$ws.Range("E51").Value = "  -0.77%  "